$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    # D-column price cells hold numeric-looking strings (e.g. "24.80",
    # "29.381.57") that must stay literal text, preserving exact digits/
    # trailing zeros. Force Text format before assignment so Excel's COM
    # layer doesn't auto-coerce the string into a Double, then restore the
    # default "Normal" style so no stray formatting is left on the cell.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

function Set-PlainValue($cellRef, $value) {
    $ws.Range($cellRef).Value = $value
}

Set-TextValue "D2" "29.381.57"
Set-PlainValue "E2" "  -0.11%  "
Set-TextValue "D3" "1.848.54"
Set-PlainValue "E3" "  -0.06%  "
Set-TextValue "D4" "0.9987"
Set-PlainValue "E4" "  -0.11%  "
Set-TextValue "D5" "240.64"
Set-PlainValue "E5" "  -0.08%  "
Set-PlainValue "E6" "  -0.49%  "
Set-TextValue "D7" "0.9997"
Set-PlainValue "E7" "  -0.08%  "
Set-PlainValue "E8" "  -1.20%  "
Set-TextValue "D9" "0.2903"
Set-PlainValue "E9" "  -1.34%  "
Set-TextValue "D10" "24.80"
Set-PlainValue "E10" "  +1.22%  "
Set-TextValue "D11" "0.07743"
Set-PlainValue "E11" "  -0.09%  "
Set-TextValue "D12" "5.029"
Set-PlainValue "E12" "  +0.07%  "
Set-TextValue "D13" "0.6801"
Set-PlainValue "E13" "  -0.02%  "
Set-TextValue "D14" "0.00001054"
Set-PlainValue "E14" "  -3.17%  "
Set-TextValue "D15" "83.04"
Set-PlainValue "E15" "  -0.81%  "
Set-TextValue "D16" "6.147"
Set-PlainValue "E16" "  -0.11%  "
Set-TextValue "D17" "29.388.65"
Set-PlainValue "E17" "  -0.16%  "
Set-TextValue "D18" "228.31"
Set-PlainValue "E18" "  -0.53%  "
Set-TextValue "D19" "12.36"
Set-PlainValue "E19" "  -0.87%  "
Set-TextValue "D20" "0.9991"
Set-PlainValue "E20" "  -0.14%  "
Set-TextValue "D21" "7.477"
Set-PlainValue "E21" "  +0.30%  "
Set-TextValue "D22" "0.9998"
Set-PlainValue "E22" "  -0.08%  "
Set-TextValue "D23" "158.95"
Set-TextValue "D24" "0.1385"
Set-PlainValue "E24" "  -0.27%  "
Set-TextValue "D25" "8.433"
Set-TextValue "D26" "17.67"
Set-PlainValue "E26" "  -0.09%  "
Set-TextValue "D27" "1.411"
Set-PlainValue "E27" "  +7.60%  "
Set-TextValue "D28" "1.458"
Set-PlainValue "E28" "  -0.65%  "
Set-TextValue "D29" "0.05610"
Set-PlainValue "E29" "  -2.55%  "
Set-TextValue "D30" "4.107"
Set-PlainValue "E30" "  -0.15%  "
Set-TextValue "D31" "4.071"
Set-PlainValue "E31" "  +0.47%  "
Set-TextValue "D32" "1.165"
Set-PlainValue "E32" "  +0.44%  "
Set-PlainValue "E33" "  -1.04%  "
Set-TextValue "D34" "0.6968"
Set-PlainValue "E34" "  -1.72%  "
Set-TextValue "D35" "2.588"
Set-PlainValue "E35" "  +0.00%  "
Set-TextValue "D36" "0.01802"
Set-PlainValue "E36" "  +0.23%  "
Set-TextValue "D37" "1.227.53"
Set-PlainValue "E37" "  -0.14%  "
Set-TextValue "D38" "2.719"
Set-PlainValue "E38" "  -2.05%  "
Set-TextValue "D39" "6.367"
Set-PlainValue "E39" "  -1.73%  "
Set-TextValue "D40" "0.9010"
Set-PlainValue "E40" "  -1.32%  "
Set-TextValue "D41" "0.9996"
Set-PlainValue "E41" "  -0.10%  "
Set-TextValue "D42" "101.40"
Set-PlainValue "E42" "  -0.33%  "
Set-TextValue "D43" "65.55"
Set-PlainValue "E43" "  -1.14%  "
Set-TextValue "D44" "7.201"
Set-PlainValue "E44" "  +0.68%  "
Set-PlainValue "E45" "  -0.48%  "
Set-TextValue "D46" "9.020"
Set-PlainValue "E46" "  -0.20%  "
Set-TextValue "D47" "1.686"
Set-PlainValue "E47" "  +0.01%  "
Set-TextValue "D48" "0.1145"
Set-PlainValue "E48" "  +1.82%  "
Set-PlainValue "B49" "Cronos"
Set-PlainValue "C49" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D49" "0.05701"
Set-PlainValue "E49" "  -0.24%  "
Set-PlainValue "B50" "BabyDogeCoin"
Set-PlainValue "C50" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D50" "0.00000000111"
Set-PlainValue "E50" "  -8.80%  "
Set-TextValue "D51" "0.4625"
Set-PlainValue "E51" "  -0.08%  "
